$wb = $excel.ActiveWorkbook

$wsAdd  = $wb.Worksheets.Item("Add Load")
$wsView = $wb.Worksheets.Item("View Load")

# Existing TC003 verification row on "View Load" now points at the updated
# Origin id. Do this first so the edited string keeps its place in the
# shared-string table (matching how Excel itself reorders strings).
$wsView.Range("G7").Value = "Alaska_1012101412"

# --- "Add Load" sheet: append two new test-data rows (TC008, TC009) ---
$wsAdd.Range("A9").Value = "Loads_TC008"
$wsAdd.Range("B9").Value = "New Day"
$wsAdd.Range("C9").Value = "Current Date"
$wsAdd.Range("D9").Value = "CP Shipper"
$wsAdd.Range("E9").Value = "TestContact"
$wsAdd.Range("F9").Value = "Corn"
$wsAdd.Range("G9").NumberFormat = "@"
$wsAdd.Range("G9").Value = "1.47"
$wsAdd.Range("H9").Value = "Bushels"
$wsAdd.Range("I9").Value = "Alaska"
$wsAdd.Range("J9").Value = "Roger"
$wsAdd.Range("K9").Value = "Added new load successfully"

$wsAdd.Range("A10").Value = "Loads_TC009"
$wsAdd.Range("B10").Value = "New Day"
$wsAdd.Range("C10").Value = "Current Date"
$wsAdd.Range("D10").Value = "CP Shipper"
$wsAdd.Range("E10").Value = "NA"
$wsAdd.Range("F10").Value = "Corn"
$wsAdd.Range("G10").NumberFormat = "@"
$wsAdd.Range("G10").Value = "1.48"
$wsAdd.Range("H10").Value = "Bushels"
$wsAdd.Range("I10").Value = "Alaska"
$wsAdd.Range("J10").Value = "Roger"
$wsAdd.Range("K10").Value = "Added new load successfully"

# --- "View Load" sheet: append the matching verification rows ---
$wsView.Range("A15").Value = "Loads_TC008"
$wsView.Range("B15").Value = "Current Date"
$wsView.Range("C15").Value = "NA"
$wsView.Range("D15").Value = "NA"
$wsView.Range("E15").Value = "New Day"
$wsView.Range("F15").Value = "NA"
$wsView.Range("G15").Value = "Alaska_1011011343"
$wsView.Range("H15").Value = "Roger"
$wsView.Range("I15").NumberFormat = "@"
$wsView.Range("I15").Value = "1.47"
$wsView.Range("J15").Value = "Bushels"
$wsView.Range("K15").Value = "Corn"
$wsView.Range("L15").Value = "NA"
$wsView.Range("M15").Value = "ADD"
$wsView.Range("N15").Value = "Webtable validated successfully"

$wsView.Range("A16").Value = "Loads_TC009"
$wsView.Range("B16").Value = "Current Date"
$wsView.Range("C16").Value = "NA"
$wsView.Range("D16").Value = "NA"
$wsView.Range("E16").Value = "New Day"
$wsView.Range("F16").Value = "NA"
$wsView.Range("G16").Value = "Alaska_1011011344"
$wsView.Range("H16").Value = "Roger"
$wsView.Range("I16").NumberFormat = "@"
$wsView.Range("I16").Value = "1.48"
$wsView.Range("J16").Value = "Bushels"
$wsView.Range("K16").Value = "Corn"
$wsView.Range("L16").Value = "NA"
$wsView.Range("M16").Value = "ADD"
$wsView.Range("N16").Value = "Webtable validated successfully"

# --- Update selections to match the authored state ---
[void]$wsAdd.Range("D9").Select()
[void]$wsView.Range("G23").Select()

# "View Load" becomes the active/visible tab (was "Add Load")
$wsView.Activate()
